# Localization status report refresh ("Generate Report for Archive").
#
# 1. The handoff status moves from "Ready for handoff" to "In Translation"
#    on every sheet (Overview's per-locale status columns, and the
#    "Status" column on each per-locale detail sheet).
# 2. The (cosmetic) width of the affected Status-ish columns is reduced
#    to better fit the new, shorter label.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Update the status values -------------------------------------------------

# Overview sheet: columns E (zh-cn) and F (de-de) hold the latest status
# for row 2 (the single data row).
if ($ws_overview.Range("E2").Value2 -eq $oldStatus) {
    $ws_overview.Range("E2").Value = $newStatus
}
if ($ws_overview.Range("F2").Value2 -eq $oldStatus) {
    $ws_overview.Range("F2").Value = $newStatus
}

# Per-locale detail sheets: column C is "Status".
if ($ws_zhcn.Range("C2").Value2 -eq $oldStatus) {
    $ws_zhcn.Range("C2").Value = $newStatus
}
if ($ws_dede.Range("C2").Value2 -eq $oldStatus) {
    $ws_dede.Range("C2").Value = $newStatus
}

# --- Shrink the columns that used to hold the longer status text --------------
# Target stored column width ~13.41 characters; 12.5 is the closest attainable
# ColumnWidth given Excel's internal pixel snapping.
$newColumnWidth = 12.5

$ws_overview.Range("E1").ColumnWidth = $newColumnWidth
$ws_overview.Range("F1").ColumnWidth = $newColumnWidth
$ws_zhcn.Range("C1").ColumnWidth = $newColumnWidth
$ws_dede.Range("C1").ColumnWidth = $newColumnWidth
